$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 10326
$ws.Range("F4").Value = 2527
$ws.Range("F8").Value = 481
$ws.Range("F9").Value = 773
$ws.Range("F12").Value = 1085
$ws.Range("F13").Value = 3218
$ws.Range("F14").Value = 2395
$ws.Range("F16").Value = 2148
$ws.Range("F20").Value = 1594
$ws.Range("F21").Value = 574
$ws.Range("F22").Value = 62
$ws.Range("F23").Value = 247
$ws.Range("F24").Value = 9
$ws.Range("F27").Value = 48
$ws.Range("F29").Value = 5
$ws.Range("F31").Value = 385
$ws.Range("F32").Value = 598
$ws.Range("F33").Value = 9
$ws.Range("F34").Value = 54
$ws.Range("F37").Value = 1576
$ws.Range("F38").Value = 461
$ws.Range("F39").Value = 443
$ws.Range("F40").Value = 1707
$ws.Range("F41").Value = 136
$ws.Range("F42").Value = 440
$ws.Range("F43").Value = 50
$ws.Range("F44").Value = 454
$ws.Range("F45").Value = 1017

# Sheet 2: 演出 (Performance)
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 46
$ws.Range("F9").Value = 7
$ws.Range("F10").Value = 1

# Sheet 4: 全部类型 (All Types)
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 10326
$ws.Range("F9").Value = 773
$ws.Range("F10").Value = 1085
$ws.Range("F11").Value = 3218
$ws.Range("F12").Value = 2395
$ws.Range("F13").Value = 2148
$ws.Range("F15").Value = 1594
$ws.Range("F16").Value = 574
$ws.Range("F17").Value = 62
$ws.Range("F18").Value = 247
$ws.Range("F19").Value = 9
$ws.Range("F22").Value = 48
$ws.Range("F24").Value = 5
$ws.Range("F26").Value = 385
$ws.Range("F27").Value = 598
$ws.Range("F28").Value = 9
$ws.Range("F29").Value = 46
$ws.Range("F32").Value = 54
$ws.Range("F35").Value = 1576
$ws.Range("F36").Value = 461
$ws.Range("F38").Value = 443
$ws.Range("F39").Value = 1707
$ws.Range("F40").Value = 136
$ws.Range("F42").Value = 7
$ws.Range("F43").Value = 1
$ws.Range("F44").Value = 440
$ws.Range("F45").Value = 50
$ws.Range("F46").Value = 454
$ws.Range("F47").Value = 1017

